# Re-order the "Recorded By" (column G) values so that the literal
# token "System" (if present in the comma-separated list) is moved to
# the front of the list, preserving the relative order of the other
# tokens.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "backup@backdoor.com, System" -> "System, backup@backdoor.com"
#          "system, backup@backdoor.com, System" -> "System, system, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value -or $value -eq "") {
        continue
    }

    $parts = $value -split ","
    $trimmedParts = @()
    foreach ($p in $parts) {
        $trimmedParts += $p.Trim()
    }

    # Case-sensitive search for the exact literal token "System".
    # NOTE: the default PowerShell comparison/contains operators in this
    # runtime are case-insensitive (and -ceq/-cne do not actually force
    # case sensitivity here), so we must use the .Equals() instance
    # method, which performs an ordinal (case-sensitive) comparison.
    $hasExactSystem = $false
    foreach ($p in $trimmedParts) {
        if ($p.Equals("System")) {
            $hasExactSystem = $true
        }
    }

    if ($hasExactSystem) {
        $removedOne = $false
        $rest = @()
        foreach ($p in $trimmedParts) {
            if ((-not $removedOne) -and $p.Equals("System")) {
                $removedOne = $true
            } else {
                $rest += $p
            }
        }
        $newParts = @("System") + $rest
        $newValue = [string]::Join(", ", $newParts)

        if (-not $newValue.Equals($value)) {
            $cell.Value2 = $newValue
        }
    }
}
